# The demo run added an "OrderNumber" header in column A and left the
# next empty row selected. This reverts/re-applies that single content
# change: a new header label in A1 (which also introduces a new shared
# string "OrderNumber"), and leaves the selection on the next blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell - becomes a new shared-string entry "OrderNumber".
$ws.Range("A1").Value = "OrderNumber"

# Reflect the final selection state (entire row 6 selected), matching
# what Excel leaves selected after entering data in the header row.
$ws.Rows(6).Select()
